$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) The "ERD:" paragraph (2nd paragraph): switch language en-US -> fr-FR
#    on the paragraph mark + the "ERD: " run, and replace the hyperlink
#    display text (diagrams.net link -> Google Drive link) while dropping
#    the per-run w:lang that used to sit on the hyperlink runs.
# -----------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$rng = $p2.Range

$newParaXml = '<w:p ' +
  'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
  'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
  'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
  'w14:paraId="6CBBF0DE" w14:textId="18204067" w:rsidR="000F3BD4" w:rsidRPr="001406DF" w:rsidRDefault="001406DF" w:rsidP="001406DF">' +
    '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="001406DF"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">ERD: </w:t></w:r>' +
    '<w:hyperlink r:id="rId5" w:history="1">' +
      '<w:r w:rsidRPr="001406DF"><w:t xml:space="preserve">https://drive.google.com/file/d/1RbrlxSOkf1P4shHWWQ5IP67-IQs8uE7h/view?usp=</w:t></w:r>' +
      '<w:r w:rsidRPr="001406DF"><w:t xml:space="preserve">s</w:t></w:r>' +
      '<w:r w:rsidRPr="001406DF"><w:t xml:space="preserve">haring</w:t></w:r>' +
    '</w:hyperlink>' +
  '</w:p>'

$rng.InsertXML($newParaXml)

# Point the hyperlink at the new Google Drive address (the object model
# always re-homes the relationship, so fetch the hyperlink fresh first).
$p2 = $d.Paragraphs.Item(2)
$hl = $d.Hyperlinks.Item(2)
$hl.Address = "https://drive.google.com/file/d/1RbrlxSOkf1P4shHWWQ5IP67-IQs8uE7h/view?usp=sharing"

# Re-apply the "Hyperlink" character style to the (now re-split) display
# text runs, since InsertXML above could not carry an rStyle across.
$hl = $d.Hyperlinks.Item(2)
$linkRange = $hl.Range
$linkRange.Style = "Hyperlink"

# -----------------------------------------------------------------------
# 2) The trailing empty paragraph: en-US -> fr-FR on its paragraph mark.
# -----------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.LanguageID = "fr-FR"
